$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "34.615.04"
$ws.Range("E2").Value2 = "  +1.25%  "

$ws.Range("D3").Value2 = "1.794.27"
$ws.Range("E3").Value2 = "  +0.72%  "

$ws.Range("E4").Value2 = "  -0.06%  "

$ws.Range("D5").Value2 = "'227.05"
$ws.Range("E5").Value2 = "  +0.57%  "

$ws.Range("D6").Value2 = "'0.558"
$ws.Range("E6").Value2 = "  +2.14%  "

$ws.Range("E7").Value2 = "  -0.10%  "

$ws.Range("D8").Value2 = "'32.99"
$ws.Range("E8").Value2 = "  +4.09%  "

$ws.Range("D9").Value2 = "'0.297"
$ws.Range("E9").Value2 = "  +2.06%  "

$ws.Range("E10").Value2 = "  +1.18%  "

$ws.Range("E11").Value2 = "  +0.40%  "

$ws.Range("D12").Value2 = "2.054.25"
$ws.Range("E12").Value2 = "  +0.72%  "

$ws.Range("D13").Value2 = "1.803.20"
$ws.Range("E13").Value2 = "  +1.19%  "

$ws.Range("D14").Value2 = "'11.06"
$ws.Range("E14").Value2 = "  +0.75%  "

$ws.Range("D15").Value2 = "'0.637"
$ws.Range("E15").Value2 = "  +2.37%  "

$ws.Range("D16").Value2 = "34.559.92"
$ws.Range("E16").Value2 = "  +1.30%  "

$ws.Range("D17").Value2 = "'4.29"
$ws.Range("E17").Value2 = "  +2.92%  "

$ws.Range("D18").Value2 = "'68.80"
$ws.Range("E18").Value2 = "  +1.32%  "

$ws.Range("D19").Value2 = "'248.30"
$ws.Range("E19").Value2 = "  +0.97%  "

$ws.Range("E20").Value2 = "  +2.08%  "

$ws.Range("D21").Value2 = "'11.27"
$ws.Range("E21").Value2 = "  +3.00%  "

$ws.Range("E22").Value2 = "  -0.13%  "

$ws.Range("E23").Value2 = "  +1.99%  "

$ws.Range("E24").Value2 = "  +1.23%  "

$ws.Range("D25").Value2 = "'165.47"
$ws.Range("E25").Value2 = "  +2.15%  "

$ws.Range("E26").Value2 = "  +1.66%  "

$ws.Range("D27").Value2 = "'16.56"
$ws.Range("E27").Value2 = "  +1.53%  "

$ws.Range("D28").Value2 = "'0.116"
$ws.Range("E28").Value2 = "  +2.32%  "

$ws.Range("D30").Value2 = "'4.12"
$ws.Range("E30").Value2 = "  +13.40%  "

$ws.Range("B31").Value2 = "PancakeSwap"
$ws.Range("C31").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value2 = "'1.24"
$ws.Range("E31").Value2 = "  +0.49%  "

$ws.Range("B32").Value2 = "Filecoin"
$ws.Range("C32").Value2 = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value2 = "'3.82"
$ws.Range("E32").Value2 = "  +2.46%  "

$ws.Range("B33").Value2 = "Hedera"
$ws.Range("C33").Value2 = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value2 = "'0.0525"
$ws.Range("E33").Value2 = "  +1.14%  "

$ws.Range("E34").Value2 = "  +2.81%  "

$ws.Range("D35").Value2 = "1.426.44"
$ws.Range("E35").Value2 = "  -1.05%  "

$ws.Range("E36").Value2 = "  +6.84%  "

$ws.Range("D37").Value2 = "'0.672"
$ws.Range("E37").Value2 = "  +2.84%  "

$ws.Range("E38").Value2 = "  +1.08%  "

$ws.Range("E39").Value2 = "  +2.15%  "

$ws.Range("D40").Value2 = "'85.41"
$ws.Range("E40").Value2 = "  +6.57%  "

$ws.Range("D41").Value2 = "'2.41"
$ws.Range("E41").Value2 = "  +1.35%  "

$ws.Range("D42").Value2 = "'0.934"
$ws.Range("E42").Value2 = "  +1.28%  "

$ws.Range("D43").Value2 = "'2.75"

$ws.Range("D44").Value2 = "'13.63"
$ws.Range("E44").Value2 = "  +0.39%  "

$ws.Range("E45").Value2 = "  +3.81%  "

$ws.Range("E46").Value2 = "  +1.25%  "

$ws.Range("E47").Value2 = "  +0.09%  "

$ws.Range("D48").Value2 = "1.953.75"
$ws.Range("E48").Value2 = "  +0.63%  "

$ws.Range("D49").Value2 = "'106.05"
$ws.Range("E49").Value2 = "  +0.61%  "

$ws.Range("B50").Value2 = "BabyDogeCoin"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value2 = "0.0₆0130"
$ws.Range("E50").Value2 = "  -5.54%  "

$ws.Range("B51").Value2 = "PaxDollar"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value2 = "'1.00"
$ws.Range("E51").Value2 = "  -0.07%  "
